$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.823791027069092
$ws.Range("B1").Value = 3.012819766998291
$ws.Range("C1").Value = 1.896028161048889
$ws.Range("D1").Value = 1.614491939544678
$ws.Range("E1").Value = 1.476139545440674
